# Add new columns I (I0) and J (IF) to the sheet, mirroring the existing
# header style used by the other header cells (B1:H1), and fill in the
# per-row numeric values for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy the formatting from an existing header cell (H1) onto the
# new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-36: column I (I0) and column J (IF)
$values = @{
    2  = @(5, 6)
    3  = @(6, 7)
    4  = @(3, 5)
    5  = @(6, 6)
    6  = @(4, 7)
    7  = @(6, 7)
    8  = @(1, 2)
    9  = @(8, 8)
    10 = @(7, 8)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(5, 7)
    14 = @(8, 8)
    15 = @(8, 9)
    16 = @(8, 9)
    17 = @(8, 8)
    18 = @(5, 6)
    19 = @(6, 7)
    20 = @(5, 6)
    21 = @(7, 7)
    22 = @(1, 3)
    23 = @(1, 6)
    24 = @(1, 5)
    25 = @(1, 5)
    26 = @(1, 4)
    27 = @(1, 4)
    28 = @(1, 5)
    29 = @(1, 2)
    30 = @(1, 4)
    31 = @(1, 5)
    32 = @(1, 5)
    33 = @(1, 4)
    34 = @(1, 4)
    35 = @(1, 2)
    36 = @(1, 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
